$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Events sheet ("sheet4"): make room for two new columns (dialogIDA /
#    dialogIDB) by re-typing the existing header/trigger values two columns
#    to the right, then filling the vacated + new cells.
# ---------------------------------------------------------------------------
$wsEvents = $wb.Worksheets.Item("Events")

# --- capture the values that need to move, before anything is overwritten ---
$row1_H = $wsEvents.Range("H1").Value()   # message
$row1_J = $wsEvents.Range("J1").Value()   # Triggers:
$row1_K = $wsEvents.Range("K1").Value()   # Attack
$row1_N = $wsEvents.Range("N1").Value()   # Harm
$row1_Q = $wsEvents.Range("Q1").Value()   # Death

$row2_K = $wsEvents.Range("K2").Value()   # IndividualID
$row2_L = $wsEvents.Range("L2").Value()   # EventID
$row2_N = $wsEvents.Range("N2").Value()   # IndividualID
$row2_O = $wsEvents.Range("O2").Value()   # EventID
$row2_Q = $wsEvents.Range("Q2").Value()   # IndividualID
$row2_R = $wsEvents.Range("R2").Value()   # EventID

$row3_H = $wsEvents.Range("H3").Value()   # "Have at you!"
$row3_K = $wsEvents.Range("K3").Value()   # 206
$row3_L = $wsEvents.Range("L3").Value()   # 2

$row4_H = $wsEvents.Range("H4").Value()   # 0

# --- clear every cell that is about to be relocated -----------------------
$wsEvents.Range("H1").Value = ""
$wsEvents.Range("J1").Value = ""
$wsEvents.Range("K1").Value = ""
$wsEvents.Range("N1").Value = ""
$wsEvents.Range("Q1").Value = ""

$wsEvents.Range("K2").Value = ""
$wsEvents.Range("L2").Value = ""
$wsEvents.Range("N2").Value = ""
$wsEvents.Range("O2").Value = ""
$wsEvents.Range("Q2").Value = ""
$wsEvents.Range("R2").Value = ""

$wsEvents.Range("H3").Value = ""
$wsEvents.Range("K3").Value = ""
$wsEvents.Range("L3").Value = ""

$wsEvents.Range("H4").Value = ""

# --- write the relocated values into their new (two-columns-right) homes --
$wsEvents.Range("J1").Value = $row1_H
$wsEvents.Range("L1").Value = $row1_J
$wsEvents.Range("M1").Value = $row1_K
$wsEvents.Range("P1").Value = $row1_N
$wsEvents.Range("S1").Value = $row1_Q

$wsEvents.Range("M2").Value = $row2_K
$wsEvents.Range("N2").Value = $row2_L
$wsEvents.Range("P2").Value = $row2_N
$wsEvents.Range("Q2").Value = $row2_O
$wsEvents.Range("S2").Value = $row2_Q
$wsEvents.Range("T2").Value = $row2_R

$wsEvents.Range("J3").Value = $row3_H
$wsEvents.Range("M3").Value = $row3_K
$wsEvents.Range("N3").Value = $row3_L

$wsEvents.Range("J4").Value = $row4_H

# --- new columns H (dialogIDA) / I (dialogIDB) -----------------------------
$wsEvents.Range("H1").Value = "dialogIDA"
$wsEvents.Range("I1").Value = "dialogIDB"

$wsEvents.Range("H3").Value = 0
$wsEvents.Range("I3").Value = 0

$wsEvents.Range("H4").Value = 0
$wsEvents.Range("I4").Value = 0

# --- sheet view: Events is no longer the active tab; selection -> I2 ------
$wsEvents.Range("I2").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Dialog sheet ("sheet1") becomes the active tab; selection -> C11:C12
# ---------------------------------------------------------------------------
$wsDialog = $wb.Worksheets.Item("Dialog")
$wsDialog.Activate() | Out-Null
$wsDialog.Range("C11:C12").Select() | Out-Null
